# Update the "Förändrad" date column (C) for rows 2-11 from serial 45179
# (2023-09-10) to serial 45180 (2023-09-11), matching the automatic update
# recorded in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    $current = $cell.Value2
    if ($current -eq 45179) {
        $cell.Value = 45180
    }
}
